# Update "想去人数" (F) and "最低票价" (G) figures across the four sheets
# of the Guangzhou comic-con info workbook, per the latest scrape run.

$wb = $excel.ActiveWorkbook

# Sheet 1: 展览 (Exhibition)
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value  = 593
$ws1.Range("F5").Value  = 747
$ws1.Range("F6").Value  = 394
$ws1.Range("G6").Value  = 48
$ws1.Range("F10").Value = 236
$ws1.Range("F11").Value = 6148
$ws1.Range("F14").Value = 513
$ws1.Range("G14").Value = 29.9
$ws1.Range("F16").Value = 551
$ws1.Range("F17").Value = 366
$ws1.Range("F21").Value = 717
$ws1.Range("F22").Value = 169
$ws1.Range("F23").Value = 99
$ws1.Range("F24").Value = 321
$ws1.Range("F25").Value = 1030
$ws1.Range("F27").Value = 1863
$ws1.Range("F28").Value = 511

# Sheet 2: 演出 (Performance)
$ws2 = $wb.Worksheets.Item("演出")
$ws2.Range("F3").Value = 274
$ws2.Range("F4").Value = 54
$ws2.Range("F5").Value = 274

# Sheet 3: 本地生活 (Local life)
$ws3 = $wb.Worksheets.Item("本地生活")
$ws3.Range("F2").Value = 270

# Sheet 4: 全部类型 (All types)
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value  = 270
$ws4.Range("F3").Value  = 593
$ws4.Range("F6").Value  = 747
$ws4.Range("F8").Value  = 394
$ws4.Range("G8").Value  = 48
$ws4.Range("F12").Value = 236
$ws4.Range("F13").Value = 6148
$ws4.Range("F16").Value = 274
$ws4.Range("F17").Value = 513
$ws4.Range("G17").Value = 29.9
$ws4.Range("F19").Value = 551
$ws4.Range("F20").Value = 366
$ws4.Range("F22").Value = 54
$ws4.Range("F25").Value = 274
$ws4.Range("F28").Value = 717
$ws4.Range("F32").Value = 169
$ws4.Range("F33").Value = 99
$ws4.Range("F34").Value = 321
$ws4.Range("F35").Value = 1030
$ws4.Range("F37").Value = 1863
$ws4.Range("F38").Value = 511

$wb.Save()
